{"js": "// Update the division-problem worksheet table: replace the text of 25\n// data cells (5 per \"problem\" row; the table also has blank spacer rows\n// that are left untouched) with their new values.\n//\n// Cells are addressed by (rowIndex, columnIndex) rather than by searching\n// for the old text, because a couple of values (\"600\u00f78=\" / \"693\u00f77=\")\n// appear more than once in the table and a blind text search-and-replace\n// would not be able to tell the occurrences apart.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// rowIndex -> [newValue for column 0, 1, 2, 3, 4]\nconst newRowValues = {\n  0: [\"581\u00f75=\", \"704\u00f75=\", \"439\u00f76=\", \"140\u00f78=\", \"383\u00f74=\"],\n  4: [\"964\u00f73=\", \"731\u00f73=\", \"855\u00f73=\", \"693\u00f77=\", \"208\u00f73=\"],\n  8: [\"734\u00f76=\", \"534\u00f75=\", \"948\u00f75=\", \"436\u00f73=\", \"670\u00f77=\"],\n  12: [\"909\u00f72=\", \"782\u00f73=\", \"553\u00f79=\", \"991\u00f76=\", \"415\u00f78=\"],\n  16: [\"434\u00f74=\", \"529\u00f79=\", \"564\u00f72=\", \"323\u00f76=\", \"880\u00f79=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newRowValues)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const values = newRowValues[rowIndex];\n  for (let colIndex = 0; colIndex < values.length; colIndex++) {\n    const cell = table.getCellOrNullObject(rowIndex, colIndex);\n    cell.value = values[colIndex];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet table: replace the text of 25\n# data cells (5 per \"problem\" row; the table also has blank spacer rows\n# that are left untouched) with their new values.\n#\n# Cells are addressed by (row, column) -- 1-based, as Word COM does --\n# rather than by searching for the old text, because a couple of values\n# (\"600\u00f78=\" / \"693\u00f77=\") appear more than once in the table and a blind\n# text search-and-replace would not be able to tell the occurrences apart.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newRowValues = @{\n    1  = @(\"581\u00f75=\", \"704\u00f75=\", \"439\u00f76=\", \"140\u00f78=\", \"383\u00f74=\")\n    5  = @(\"964\u00f73=\", \"731\u00f73=\", \"855\u00f73=\", \"693\u00f77=\", \"208\u00f73=\")\n    9  = @(\"734\u00f76=\", \"534\u00f75=\", \"948\u00f75=\", \"436\u00f73=\", \"670\u00f77=\")\n    13 = @(\"909\u00f72=\", \"782\u00f73=\", \"553\u00f79=\", \"991\u00f76=\", \"415\u00f78=\")\n    17 = @(\"434\u00f74=\", \"529\u00f79=\", \"564\u00f72=\", \"323\u00f76=\", \"880\u00f79=\")\n}\n\nforeach ($rowIndex in $newRowValues.Keys) {\n    $values = $newRowValues[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n"}
